# This script updates the metrics worksheet so that every model row (rows 2-26)
# now reports the same B:Q metric values (the new "lm" training run converged to
# identical results for every ensemble member), while the model names in column A
# are re-shuffled to reflect the new training order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model names (column A) for rows 2..26, in order.
$modelNames = @(
    "model_32_9_0",
    "model_32_9_22",
    "model_32_9_21",
    "model_32_9_20",
    "model_32_9_19",
    "model_32_9_18",
    "model_32_9_17",
    "model_32_9_16",
    "model_32_9_15",
    "model_32_9_14",
    "model_32_9_13",
    "model_32_9_23",
    "model_32_9_12",
    "model_32_9_10",
    "model_32_9_9",
    "model_32_9_8",
    "model_32_9_7",
    "model_32_9_6",
    "model_32_9_5",
    "model_32_9_4",
    "model_32_9_3",
    "model_32_9_2",
    "model_32_9_1",
    "model_32_9_11",
    "model_32_9_24"
)

# New metric values (columns B..Q), identical for every row.
# (written in plain decimal form - this runtime's PowerShell parser does not
# accept scientific-notation numeric literals like 1.23e-05)
$metricValues = @(
    0.9999949039827802,
    0.9991177146462858,
    0.9999826080166279,
    0.9999986869145568,
    0.9999896941572535,
    0.000004756908385316126,
    0.0008235746498410597,
    0.00002527387479370727,
    0.0000008309921717876705,
    0.00001305243338010575,
    0.00008583594873526996,
    0.002181033788210565,
    1.000004892176531,
    0.00227388491193704,
    122.5118251992398,
    182.2367406177816
)

for ($i = 0; $i -lt $modelNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $modelNames[$i]
    for ($c = 0; $c -lt $metricValues.Length; $c++) {
        $ws.Cells.Item($row, $c + 2).Value = $metricValues[$c]
    }
}
